$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 512.4
$ws.Range("I20").Value = 512.4
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 512.4
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -282.4
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 512.4
$ws.Range("I35").Value = 512.4
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 512.4
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -133.4
$ws.Range("N35").ClearContents()
$ws.Range("H43").Value = 4179.2
$ws.Range("J43").Value = 2947
$ws.Range("L43").Value = 2947
$ws.Range("N43").Value = -3085
$ws.Range("H45").Value = 1777.5
$ws.Range("I45").Value = 1555
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 4665
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -4473
$ws.Range("N45").Value = -6384
$ws.Range("H55").Value = 229.77777
$ws.Range("J55").Value = 256.66666
$ws.Range("L55").Value = 256.66666
$ws.Range("N55").Value = -684.66666
$ws.Range("H92").Value = 2541.5
$ws.Range("I92").Value = 2541.5
$ws.Range("K92").Value = 2541.5
$ws.Range("M92").Value = -1293.5
$ws.Range("H98").Value = 813.55554
$ws.Range("I98").Value = 813.55554
$ws.Range("K98").Value = 813.55554
$ws.Range("M98").Value = 684.44446
$ws.Range("H106").Value = 2299.5
$ws.Range("I106").Value = 2299.5
$ws.Range("K106").Value = 2299.5
$ws.Range("M106").Value = -1668.5
$ws.Range("H111").Value = 800
$ws.Range("I111").Value = 200
$ws.Range("J111").Value = 1100
$ws.Range("K111").Value = 600
$ws.Range("L111").Value = 3300
$ws.Range("M111").Value = 2467
$ws.Range("N111").Value = -9434
$ws.Range("H122").Value = 813.55554
$ws.Range("I122").Value = 813.55554
$ws.Range("K122").Value = 2440.66662
$ws.Range("M122").Value = 9.333380000000034
$ws.Range("H132").Value = 2574.8333
$ws.Range("I132").Value = 2167.6428
$ws.Range("K132").Value = 6502.928400000001
$ws.Range("M132").Value = -3972.928400000001
$ws.Range("H138").Value = 2669.889
$ws.Range("J138").Value = 2989.4736
$ws.Range("L138").Value = 8968.4208
$ws.Range("N138").Value = -19248.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 998.6667
$ws.Range("I122").Value = 998.6667
$ws.Range("K122").Value = 2996.0001
$ws.Range("M122").Value = -546.0001000000002
$ws.Range("H132").Value = 4151.4287
$ws.Range("I132").Value = 4206.6
$ws.Range("J132").Value = 4013.5
$ws.Range("K132").Value = 12619.8
$ws.Range("L132").Value = 12040.5
$ws.Range("M132").Value = -10089.8
$ws.Range("N132").Value = -17100.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4516.3335
$ws.Range("I86").Value = 9999
$ws.Range("J86").Value = 1775
$ws.Range("K86").Value = 9999
$ws.Range("L86").Value = 1775
$ws.Range("M86").Value = -8876
$ws.Range("N86").Value = -4021
$ws.Range("H89").Value = 4516.3335
$ws.Range("I89").Value = 9999
$ws.Range("J89").Value = 1775
$ws.Range("K89").Value = 49995
$ws.Range("L89").Value = 8875
$ws.Range("M89").Value = -44379
$ws.Range("N89").Value = -20107
$ws.Range("H105").Value = 4360
$ws.Range("I105").Value = 3433.3333
$ws.Range("J105").Value = 5750
$ws.Range("K105").Value = 3433.3333
$ws.Range("L105").Value = 5750
$ws.Range("M105").Value = -1686.3333
$ws.Range("N105").Value = -9244
$ws.Range("H107").Value = 1150
$ws.Range("I107").Value = 1116.6666
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1116.6666
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 803.3334
$ws.Range("N107").Value = -5040
$ws.Range("H134").Value = 2465.125
$ws.Range("I134").Value = 2569.4666
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 7708.399800000001
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -5173.399800000001
$ws.Range("N134").Value = -7770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1315.6666
$ws.Range("I58").Value = 750
$ws.Range("J58").Value = 1598.5
$ws.Range("K58").Value = 750
$ws.Range("L58").Value = 1598.5
$ws.Range("M58").Value = -547
$ws.Range("N58").Value = -2004.5
$ws.Range("H86").Value = 668233
$ws.Range("I86").Value = 561879.8
$ws.Range("K86").Value = 561879.8
$ws.Range("M86").Value = -560756.8
$ws.Range("H89").Value = 668233
$ws.Range("I89").Value = 561879.8
$ws.Range("K89").Value = 2809399
$ws.Range("M89").Value = -2803783
$ws.Range("H105").Value = 11838
$ws.Range("I105").Value = 13067.75
$ws.Range("K105").Value = 13067.75
$ws.Range("M105").Value = -11320.75
$ws.Range("H134").Value = 2080.9375
$ws.Range("I134").Value = 2018.7333
$ws.Range("K134").Value = 6056.199900000001
$ws.Range("M134").Value = -3521.199900000001
$ws.Range("H136").Value = 1315.6666
$ws.Range("I136").Value = 750
$ws.Range("J136").Value = 1598.5
$ws.Range("K136").Value = 2250
$ws.Range("L136").Value = 4795.5
$ws.Range("M136").Value = 300
$ws.Range("N136").Value = -9895.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 835
$ws.Range("I92").Value = 587.5
$ws.Range("K92").Value = 1762.5
$ws.Range("M92").Value = -514.5
$ws.Range("H98").Value = 512.25
$ws.Range("I98").Value = 450
$ws.Range("J98").Value = 533
$ws.Range("K98").Value = 1350
$ws.Range("L98").Value = 1599
$ws.Range("M98").Value = 148
$ws.Range("N98").Value = -4595
$ws.Range("H131").Value = 1412.4375
$ws.Range("J131").Value = 1412.4375
$ws.Range("L131").Value = 4237.3125
$ws.Range("N131").Value = -14317.3125
$ws.Range("H137").Value = 2382.2
$ws.Range("I137").Value = 1310.8
$ws.Range("K137").Value = 3932.4
$ws.Range("M137").Value = 1167.6
$ws.Range("H139").Value = 9995
$ws.Range("I139").Value = 9995
$ws.Range("K139").Value = 29985
$ws.Range("M139").Value = -24845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2002.5
$ws.Range("I102").Value = 2002.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2002.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -380.5
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15470
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2387.5
$ws.Range("I40").Value = 2387.5
$ws.Range("K40").Value = 2387.5
$ws.Range("M40").Value = -2251.5
$ws.Range("H46").Value = 4999.933
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3812
$ws.Range("H96").Value = 348000
$ws.Range("J96").Value = 348000
$ws.Range("L96").Value = 348000
$ws.Range("N96").Value = -353492
$ws.Range("H122").Value = 8639.8
$ws.Range("I122").Value = 8299.75
$ws.Range("K122").Value = 24899.25
$ws.Range("M122").Value = -22449.25
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920
$ws.Range("H132").Value = 6899.5
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 1546.3846
$ws.Range("I136").Value = 1546.3846
$ws.Range("K136").Value = 4639.1538
$ws.Range("M136").Value = -2089.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 23722
$ws.Range("I75").Value = 20000
$ws.Range("J75").Value = 24962.666
$ws.Range("K75").Value = 20000
$ws.Range("L75").Value = 24962.666
$ws.Range("M75").Value = -19064
$ws.Range("N75").Value = -26834.666
$ws.Range("H78").Value = 23722
$ws.Range("I78").Value = 20000
$ws.Range("J78").Value = 24962.666
$ws.Range("K78").Value = 60000
$ws.Range("L78").Value = 74887.998
$ws.Range("M78").Value = -55320
$ws.Range("N78").Value = -84247.998
$ws.Range("H105").Value = 67795
$ws.Range("J105").Value = 67795
$ws.Range("L105").Value = 67795
$ws.Range("N105").Value = -74783
$ws.Range("H122").Value = 1100
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 2476.2222
$ws.Range("I126").Value = 1831.5
$ws.Range("J126").Value = 3765.6667
$ws.Range("K126").Value = 5494.5
$ws.Range("L126").Value = 11297.0001
$ws.Range("M126").Value = -3024.5
$ws.Range("N126").Value = -16237.0001
$ws.Range("H132").Value = 1008.6923
$ws.Range("I132").Value = 1008.6923
$ws.Range("K132").Value = 3026.0769
$ws.Range("M132").Value = -496.0769
